$d = $word.ActiveDocument

# 1. Replace the title of the first film.
$d.Content.Find.Execute("Vingadores", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Missão Impossível", 2)

# 2. Remove the middle paragraphs (Superman, Batman, Starwars) entirely,
#    keeping the last one's paragraph mark so the document still ends
#    with a trailing empty paragraph.
$pFirstToRemove = $d.Paragraphs.Item(2)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($pFirstToRemove.Range.Start, $pLast.Range.Start).Delete()

# 3. Clear the text of what is now the second (last) paragraph
#    ("Harry Potter"), leaving it as an empty paragraph.
$pLast2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($pLast2.Range.Start, $pLast2.Range.End - 1).Delete()
